$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.495.61"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").Value = "1.914.07"
$ws.Range("D4").Value = "'1.0000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'244.76"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").Value = "'0.4842"
$ws.Range("E7").Value = "  +3.64%  "
$ws.Range("D8").Value = "'0.2890"
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("D9").Value = "'0.06709"
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").Value = "'109.68"
$ws.Range("E10").Value = "  +2.75%  "
$ws.Range("D11").Value = "'18.92"
$ws.Range("E11").Value = "  +5.49%  "
$ws.Range("D12").Value = "1.915.42"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").Value = "'0.07551"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("E15").Value = "  +2.54%  "
$ws.Range("D16").Value = "'275.87"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("D17").Value = "30.506.19"
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "'0.000007536"
$ws.Range("E19").Value = "  -0.44%  "
$ws.Range("E20").Value = "  -0.71%  "
$ws.Range("D21").Value = "2.159.80"
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").Value = "'5.480"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").Value = "'6.441"
$ws.Range("E24").Value = "  +4.27%  "
$ws.Range("D25").Value = "'9.433"
$ws.Range("D26").Value = "'163.68"
$ws.Range("E26").Value = "  -2.62%  "
$ws.Range("D27").Value = "'20.12"
$ws.Range("E27").Value = "  -5.76%  "
$ws.Range("D28").Value = "'2.123"
$ws.Range("E28").Value = "  +4.80%  "
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").Value = "'1.399"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").Value = "'4.148"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").Value = "'4.049"
$ws.Range("E32").Value = "  +2.98%  "
$ws.Range("D33").Value = "'0.04991"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("D34").Value = "'0.7291"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").Value = "'1.130"
$ws.Range("E35").Value = "  -0.90%  "
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("D37").Value = "'2.734"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").Value = "'110.76"
$ws.Range("E40").Value = "  +1.85%  "
$ws.Range("D41").Value = "'2.017"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("D42").Value = "'0.4424"
$ws.Range("E42").Value = "  +5.79%  "
$ws.Range("D43").Value = "'0.8658"
$ws.Range("E43").Value = "  -0.68%  "
$ws.Range("D44").Value = "'5.834"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'67.74"
$ws.Range("E46").Value = "  +0.90%  "
$ws.Range("D47").Value = "'7.384"
$ws.Range("E47").Value = "  +3.93%  "
$ws.Range("D48").Value = "'9.274"
$ws.Range("E48").Value = "  +1.82%  "
$ws.Range("D49").Value = "'0.1239"
$ws.Range("E49").Value = "  +3.06%  "
$ws.Range("D50").Value = "'47.79"
$ws.Range("E50").Value = "  -9.65%  "
$ws.Range("D51").Value = "'1.463"
$ws.Range("E51").Value = "  +6.99%  "
